$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.518.98"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.573.07"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.95"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3732"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.39"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3332"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.135"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07482"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.93"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.987"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.923"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "1.572.78"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001119"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.36"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06788"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.394"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.49"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "22.458.59"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.390"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.575"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.71"
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.74"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.25"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "1.751.21"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.057"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.170"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.014"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.683"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08312"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02458"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2269"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06392"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.398"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.291"
$ws.Range("E41").Value = "  -4.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.33"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6299"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6156"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.050"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.42"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.217"
$ws.Range("E50").Value = "  -2.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07280"
$ws.Range("E51").Value = "  -0.29%  "
